$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "30.262.28"
Set-TextValue "E2" "  +0.47%  "
Set-TextValue "D3" "1.864.95"
Set-TextValue "E3" "  +0.27%  "
Set-TextValue "E4" "  +0.10%  "
Set-TextValue "D5" "237.32"
Set-TextValue "E5" "  +1.65%  "
Set-TextValue "D6" "1.001"
Set-TextValue "E6" "  +0.10%  "
Set-TextValue "D7" "0.4685"
Set-TextValue "E7" "  +0.58%  "
Set-TextValue "D8" "0.2865"
Set-TextValue "E8" "  +2.03%  "
Set-TextValue "D9" "0.06547"
Set-TextValue "E9" "  +0.04%  "
Set-TextValue "D10" "22.21"
Set-TextValue "E10" "  +13.28%  "
Set-TextValue "D11" "0.07897"
Set-TextValue "E11" "  +1.26%  "
Set-TextValue "D12" "97.79"
Set-TextValue "E12" "  +1.39%  "
Set-TextValue "D13" "1.867.61"
Set-TextValue "E13" "  +0.43%  "
Set-TextValue "D14" "5.173"
Set-TextValue "E14" "  +0.85%  "
Set-TextValue "D15" "0.6808"
Set-TextValue "E15" "  +2.33%  "
Set-TextValue "D16" "278.55"
Set-TextValue "E16" "  -0.75%  "
Set-TextValue "D17" "30.269.64"
Set-TextValue "E17" "  +0.39%  "
Set-TextValue "D18" "13.61"
Set-TextValue "E18" "  +8.16%  "
Set-TextValue "E19" "  +0.03%  "
Set-TextValue "B20" "ShibaInu"
Set-TextValue "C20" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D20" "0.000007350"
Set-TextValue "E20" "  +1.75%  "
Set-TextValue "B21" "Uniswap"
Set-TextValue "C21" "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue "D21" "5.391"
Set-TextValue "E21" "  -2.14%  "
Set-TextValue "D22" "2.114.52"
Set-TextValue "E22" "  +0.76%  "
Set-TextValue "E23" "  +0.01%  "
Set-TextValue "D24" "6.195"
Set-TextValue "E24" "  +1.28%  "
Set-TextValue "D25" "168.73"
Set-TextValue "E25" "  +1.89%  "
Set-TextValue "D26" "9.292"
Set-TextValue "E26" "  -0.22%  "
Set-TextValue "D27" "19.12"
Set-TextValue "E27" "  +1.49%  "
Set-TextValue "D28" "1.943"
Set-TextValue "E28" "  +1.75%  "
Set-TextValue "D29" "1.385"
Set-TextValue "D30" "0.09841"
Set-TextValue "E30" "  +3.28%  "
Set-TextValue "D31" "4.395"
Set-TextValue "E31" "  -0.06%  "
Set-TextValue "D32" "1.480"
Set-TextValue "E32" "  +0.90%  "
Set-TextValue "D33" "4.070"
Set-TextValue "E33" "  -0.66%  "
Set-TextValue "D34" "0.04747"
Set-TextValue "E34" "  +2.11%  "
Set-TextValue "D35" "1.142"
Set-TextValue "E35" "  +4.75%  "
Set-TextValue "D36" "0.7086"
Set-TextValue "E36" "  +1.12%  "
Set-TextValue "D37" "2.709"
Set-TextValue "E37" "  +0.22%  "
Set-TextValue "D38" "0.01878"
Set-TextValue "E38" "  +1.22%  "
Set-TextValue "D39" "2.623"
Set-TextValue "E39" "  +4.64%  "
Set-TextValue "D40" "76.87"
Set-TextValue "E40" "  +5.20%  "
Set-TextValue "D41" "6.305"
Set-TextValue "E41" "  +0.31%  "
Set-TextValue "D42" "1.962"
Set-TextValue "E42" "  +2.16%  "
Set-TextValue "D43" "0.8509"
Set-TextValue "E43" "  -0.09%  "
Set-TextValue "D44" "0.4189"
Set-TextValue "E44" "  +1.00%  "
Set-TextValue "D45" "1.0000"
Set-TextValue "E45" "  +0.07%  "
Set-TextValue "D46" "103.32"
Set-TextValue "E46" "  -0.30%  "
Set-TextValue "B47" "Aptos"
Set-TextValue "C47" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.221"
Set-TextValue "E47" "  +0.62%  "
Set-TextValue "B48" "Maker"
Set-TextValue "C48" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D48" "956.86"
Set-TextValue "E48" "  -3.74%  "
Set-TextValue "D49" "9.291"
Set-TextValue "E49" "  +0.52%  "
Set-TextValue "D50" "34.27"
Set-TextValue "E50" "  +0.56%  "
Set-TextValue "D51" "0.05641"
